$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kayitlar")

# Shift the existing record (old row 2) down to row 3, preserving it as text
$ws.Range("A3").Value = "'2"
$ws.Range("B3").Value = "'2025-07-16"
$ws.Range("C3").Value = "Merkez"
$ws.Range("D3").Value = "'3"
$ws.Range("E3").Value = "'2"
$ws.Range("F3").Value = "Cins D."
$ws.Range("G3").Value = "Göktan ELGÜL"
$ws.Range("A3:G3").Style = "Normal"

# New record goes into row 2
$ws.Range("A2").Value = "'3"
$ws.Range("B2").Value = "'2025-07-16"
$ws.Range("C2").Value = "İlçe"
$ws.Range("D2").Value = "'2"
$ws.Range("E2").Value = "'2"
$ws.Range("F2").Value = "Tevhid"
$ws.Range("G2").Value = "Gökhan ELGÜL"
$ws.Range("A2:G2").Style = "Normal"
